# Apply updated crypto price/volume figures as exact text values (matching the
# original inline-string cell contents), since Excel would otherwise reinterpret
# numeric-looking strings (e.g. "415.17") as actual numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "62.067.69"
Set-TextCell $ws.Range("E2") "  +8.57%  "

Set-TextCell $ws.Range("D3") "3.450.54"
Set-TextCell $ws.Range("E3") "  +5.89%  "

Set-TextCell $ws.Range("E4") "  -0.03%  "

Set-TextCell $ws.Range("D5") "415.17"
Set-TextCell $ws.Range("E5") "  +4.33%  "

Set-TextCell $ws.Range("D6") "124.44"
Set-TextCell $ws.Range("E6") "  +14.17%  "

Set-TextCell $ws.Range("D7") "3.442.44"
Set-TextCell $ws.Range("E7") "  +5.79%  "

Set-TextCell $ws.Range("D8") "0.594"
Set-TextCell $ws.Range("E8") "  +2.33%  "

Set-TextCell $ws.Range("E9") "  +0.02%  "

Set-TextCell $ws.Range("D10") "0.659"
Set-TextCell $ws.Range("E10") "  +6.19%  "

Set-TextCell $ws.Range("E11") "  +32.32%  "

Set-TextCell $ws.Range("D12") "41.55"
Set-TextCell $ws.Range("E12") "  +5.51%  "

Set-TextCell $ws.Range("E13") "  -0.28%  "

Set-TextCell $ws.Range("D14") "3.977.69"
Set-TextCell $ws.Range("E14") "  +5.42%  "

Set-TextCell $ws.Range("D15") "8.53"
Set-TextCell $ws.Range("E15") "  +2.99%  "

Set-TextCell $ws.Range("E16") "  +4.42%  "

Set-TextCell $ws.Range("D17") "3.446.09"
Set-TextCell $ws.Range("E17") "  +5.57%  "

Set-TextCell $ws.Range("D18") "61.940.71"
Set-TextCell $ws.Range("E18") "  +8.70%  "

Set-TextCell $ws.Range("E19") "  -0.03%  "

Set-TextCell $ws.Range("D20") "10.99"
Set-TextCell $ws.Range("E20") "  -0.84%  "

Set-TextCell $ws.Range("D21") "0.0000130"
Set-TextCell $ws.Range("E21") "  +19.76%  "

Set-TextCell $ws.Range("E22") "  +0.62%  "

Set-TextCell $ws.Range("D23") "82.36"
Set-TextCell $ws.Range("E23") "  +10.93%  "

Set-TextCell $ws.Range("D24") "315.23"
Set-TextCell $ws.Range("E24") "  +7.03%  "

Set-TextCell $ws.Range("D25") "13.01"
Set-TextCell $ws.Range("E25") "  +0.28%  "

Set-TextCell $ws.Range("D26") "3.18"
Set-TextCell $ws.Range("E26") "  -0.25%  "

Set-TextCell $ws.Range("D27") "31.08"
Set-TextCell $ws.Range("E27") "  +10.47%  "

Set-TextCell $ws.Range("D28") "7.85"
Set-TextCell $ws.Range("E28") "  +5.71%  "

Set-TextCell $ws.Range("D29") "7.89"
Set-TextCell $ws.Range("E29") "  -0.49%  "

Set-TextCell $ws.Range("E30") "  -2.11%  "

Set-TextCell $ws.Range("D31") "0.174"
Set-TextCell $ws.Range("E31") "  +2.64%  "

Set-TextCell $ws.Range("E32") "  +4.60%  "

Set-TextCell $ws.Range("D33") "11.59"
Set-TextCell $ws.Range("E33") "  +3.41%  "

Set-TextCell $ws.Range("E34") "  +20.58%  "

Set-TextCell $ws.Range("D35") "42.40"
Set-TextCell $ws.Range("E35") "  +5.45%  "

Set-TextCell $ws.Range("D36") "1.00"
Set-TextCell $ws.Range("E36") "  +0.00%  "

Set-TextCell $ws.Range("E37") "  -0.94%  "

Set-TextCell $ws.Range("D38") "52.46"
Set-TextCell $ws.Range("E38") "  +2.21%  "

Set-TextCell $ws.Range("D39") "3.53"
Set-TextCell $ws.Range("E39") "  +1.74%  "

Set-TextCell $ws.Range("D40") "0.997"
Set-TextCell $ws.Range("E40") "  -0.32%  "

Set-TextCell $ws.Range("E41") "  -0.08%  "

Set-TextCell $ws.Range("D42") "2.00"
Set-TextCell $ws.Range("E42") "  +6.88%  "

Set-TextCell $ws.Range("E43") "  +3.56%  "

Set-TextCell $ws.Range("D44") "134.33"
Set-TextCell $ws.Range("E44") "  -1.58%  "

Set-TextCell $ws.Range("D45") "17.31"
Set-TextCell $ws.Range("E45") "  +2.90%  "

Set-TextCell $ws.Range("D46") "0.285"
Set-TextCell $ws.Range("E46") "  +0.21%  "

Set-TextCell $ws.Range("D47") "3.90"
Set-TextCell $ws.Range("E47") "  -1.38%  "

Set-TextCell $ws.Range("D48") "22.13"
Set-TextCell $ws.Range("E48") "  -1.47%  "

Set-TextCell $ws.Range("E49") "  -0.60%  "

Set-TextCell $ws.Range("D50") "2.206.95"

Set-TextCell $ws.Range("D51") "3.783.76"
Set-TextCell $ws.Range("E51") "  +5.57%  "
